# =========================================================================
# Adds "Sheet2" (a substitution-cipher / binary-encoding worksheet) right
# after Sheet1, populates it with the cryptanalysis scratch-work, and
# updates the view state so Sheet2 ends up the active/selected sheet
# (matching the author's "stopped here for the night" checkpoint).
#
# Cell values are entered in the same order the author appears to have
# typed them (reconstructed from the shared-string table order in the
# target file) so the regenerated workbook lines up as closely as
# possible with the real edit session.
# =========================================================================

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# --- Create Sheet2 right after Sheet1 -----------------------------------
$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "Sheet2"

# --- Cryptext hex dumps (column C, rows 2-3) ------------------------------
$ws.Range("C2").Value = "0809030206071A171A081C07141D"
$ws.Range("C3").Value = "00051311191119070D091B08130B"

# --- Labels typed next -----------------------------------------------------
$ws.Range("D1").Value = "length"
$ws.Range("F1").Value = "In binary:"

# --- Candidate plaintext messages (column C, rows 5-13) -------------------
$ws.Range("C5").Value  = "GORGEOUS SUSAN"
$ws.Range("C6").Value  = "NICOTINE IS BAD"
$ws.Range("C7").Value  = "MARIJUANAS LEGAL"
$ws.Range("C8").Value  = "JUSTINE TRUDEAU"
$ws.Range("C9").Value  = "FLOYD MAYWEATHER"
$ws.Range("C10").Value = "ANGELINA JOLIE"
$ws.Range("C11").Value = "EMBEZZLED FUNDS"
$ws.Range("C12").Value = "NANETTE WORKMAN"
$ws.Range("C13").Value = "ELIZABETH MAY"

# --- Row labels for the two cryptext rows ---------------------------------
$ws.Range("B2").Value = "Cryptext 1"
$ws.Range("B3").Value = "Cryptext 2"

# --- Two more candidate messages -------------------------------------------
$ws.Range("C14").Value = "GRANT US PEACE"
$ws.Range("C15").Value = "WE'RE AWESOME"

$ws.Range("B5").Value = "Potential Messages"

# --- Binary renderings of the two cryptexts (column F) --------------------
# These are long strings of 0/1 digits -- force Text format first so Excel
# doesn't coerce them into (lossy) scientific-notation numbers.
$ws.Range("F2:F15").NumberFormat = "@"
$ws.Range("F2").Value = "0000100000001001000000110000001000000110000001110001101000010111000110100000100000011100000001110001010000011101"
$ws.Range("F3").Value = "0000000000000101000100110001000100011001000100010001100100000111000011010000100100011011000010000001001100001011"

$ws.Range("G1").Value = "Length"

# --- Binary-encoded candidate messages (column I) --------------------------
$ws.Range("I5").Value  = "01001111 01000110 01010001 01000101 01000011 01001000 01001111 01000100 00111010 01011011 01001001 01010100 01010101 01010011 00000010"
$ws.Range("I13").Value = "01001101 01000101 01001010 01011000 01000111 01000101 01011111 01000011 01010010 00101000 01010001 01000110 01001101 00010111 00000010"
$ws.Range("I6").Value  = "01000110 01000000 01000000 01001101 01010010 01001110 01010100 01010010 00111010 01000001 01001111 00100111 01010110 01011100 01001100 00000011 00001001"
$ws.Range("I7").Value  = "01000101 01001000 01010001 01001011 01001100 01010010 01011011 01011001 01011011 01011011 00111100 01001011 01010001 01011010 01001001 01000101 00001001 00001000"
$ws.Range("I8").Value  = "01000010 01011100 01010000 01010110 01001111 01001001 01011111 00110111 01001110 01011010 01001001 01000011 01010001 01011100 01011101"
$ws.Range("I9").Value  = "01001110 01000101 01001100 01011011 01000010 00100111 01010111 01010110 01000011 01011111 01011001 01000110 01000000 01010101 01001101 01011011"
$ws.Range("I10").Value = "01001001 01000111 01000100 01000111 01001010 01001110 01010100 01010110 00111010 01000010 01010011 01001011 01011101 01011000"
$ws.Range("I11").Value = "01001101 01000100 01000001 01000111 01011100 01011101 01010110 01010010 01011110 00101000 01011010 01010010 01011010 01011001 01011011"
$ws.Range("I12").Value = "01000110 01001000 01001101 01000111 01010010 01010011 01011111 00110111 01001101 01000111 01001110 01001100 01011001 01011100 01000110"
$ws.Range("I14").Value = "01001111 01011011 01000010 01001100 01010010 00100111 01001111 01000100 00111010 01011000 01011001 01000110 01010111 01011000"
$ws.Range("I15").Value = "01011111 01001100 00100100 01010000 01000011 00100111 01011011 01000000 01011111 01011011 01010011 01001010 01010001"

$ws.Range("I1").Value = "Corresponding Key for M1"
$ws.Range("L1").Value = "Decrypted M2:"

$ws.Range("B11").Value = "Found it!"

# --- Formulas (LEN helper columns) -----------------------------------------
# Row 2/3 are standalone formulas; rows 3-15 (G, J) and 5-15 (D) are filled
# as a block so the engine groups them into shared formulas, matching how
# Excel records a drag-fill down a column.
$ws.Range("D2").Formula = "=LEN(C2)"
$ws.Range("D3").Formula = "=LEN(C3)"
$ws.Range("G2").Formula = "=LEN(F2)"
$ws.Range("J2").Formula = "=LEN(I2)"
$ws.Range("G3:G15").Formula = "=LEN(F3)"
$ws.Range("J3:J15").Formula = "=LEN(I3)"
$ws.Range("D5:D15").Formula = "=LEN(C5)"

# --- Column widths (best-fit, as captured by the author) ------------------
$ws.Columns.Item("B").ColumnWidth = 18.42578125
$ws.Columns.Item("C").ColumnWidth = 30.7109375
$ws.Columns.Item("D").ColumnWidth = 9
$ws.Columns.Item("G").ColumnWidth = 8.42578125
$ws.Columns.Item("I").ColumnWidth = 24.7109375

# --- View state: Sheet2 becomes the active sheet/selection ----------------
$ws.Activate()
$ws.Range("I11").Select()
